# Add a new work-log entry row documenting the English version of the
# database work ("Tietokanta englanniksi"), continuing the existing table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "Tietokanta"
$ws.Range("D6").Value = "Tietokanta englanniksi"
$ws.Range("E6").Value = 0.5

# Leave the cursor a couple of rows below the newly entered data, as the
# author would after finishing data entry.
$ws.Range("D9").Select() | Out-Null
